$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Ajo" (garlic) at the
# "Feria Lagunitas de Puerto Montt" market. It belongs at row 281
# (most recent first), pushing the existing rows 281-287 down to 282-288.
$ws.Rows.Item(281).Insert()

$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 44757
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100112003
$ws.Range("G281").Value = "Ajo"
$ws.Range("H281").Value = "Chino"
$ws.Range("I281").Value = "Primera"
$ws.Range("J281").Value = 120
$ws.Range("K281").Value = 28000
$ws.Range("L281").Value = 30000
$ws.Range("M281").Value = 29000
$ws.Range("N281").Value = "$/caja 10 kilos"
$ws.Range("O281").Value = "China"
$ws.Range("P281").Value = 2900
$ws.Range("Q281").Value = 10
$ws.Range("R281").Value = "Hortaliza"
